# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F (same updates apply to both sheets)
$updates = @{
    2  = 1067
    3  = 352
    4  = 1468
    5  = 8681
    6  = 86
    11 = 7
    12 = 3526
    14 = 360
    15 = 73
    16 = 1124
    18 = 1109
    20 = 193
    21 = 2252
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
